# Add a new entry to the "books" table:
#   1*, Aditya Dahiya,
#   "Solutions Manual (and Beyond) for ggplot2: Elegant Graphics for Data Analysis (3e)",
#   https://aditya-dahiya.github.io/ggplot2book3e/, In progress
# This becomes the new row 4 (right after the existing "1*" ggplot2-book-solutions row),
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 4 (and everything below it) down by one row, carrying values/formats along.
[void]$ws.Rows("4:4").Insert()

# Fill the newly freed row 4 with the new book entry.
$ws.Range("A4").Value = "1*"
$ws.Range("B4").Value = "Aditya Dahiya"
$ws.Range("C4").Value = "Solutions Manual (and Beyond) for ggplot2: Elegant Graphics for Data Analysis (3e)"
$ws.Range("D4").Value = "https://aditya-dahiya.github.io/ggplot2book3e/"
$ws.Range("E4").Value = "In progress"

# Match the taller row used by the other wrapped-text entries.
$ws.Rows("4:4").RowHeight = 51

# Reflect the author's selection after the edit.
[void]$ws.Range("A4:E4").Select()

Write-Host "Inserted new ggplot2 solutions-manual row at row 4"
